$d = $word.ActiveDocument

$pairs = @(
    @("2023-12-08 Friday", "2023-12-09 Saturday"),
    @("81×58=4698", "89×87=7743"),
    @("72×28=2016", "20×86=1720"),
    @("97×12=1164", "69×94=6486"),
    @("62×17=1054", "52×31=1612"),
    @("15×49=735", "62×56=3472"),
    @("57×47=2679", "60×76=4560"),
    @("16×12=192", "26×53=1378"),
    @("50×43=2150", "11×98=1078"),
    @("84×84=7056", "45×80=3600"),
    @("92×40=3680", "96×26=2496"),
    @("48×24=1152", "44×84=3696"),
    @("75×19=1425", "70×56=3920"),
    @("48×31=1488", "76×77=5852"),
    @("37×65=2405", "16×45=720"),
    @("49×82=4018", "48×36=1728"),
    @("94×67=6298", "71×58=4118"),
    @("94×80=7520", "71×26=1846"),
    @("87×15=1305", "97×15=1455"),
    @("61×54=3294", "33×85=2805"),
    @("77×88=6776", "84×38=3192"),
    @("94×91=8554", "35×35=1225"),
    @("93×29=2697", "44×75=3300"),
    @("74×15=1110", "27×13=351"),
    @("57×48=2736", "70×86=6020"),
    @("42×98=4116", "39×35=1365")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, `
                         $true, 1, $false, $new, 2)
}

$d.Save()
